# This script applies the "Updated symbol list" commit to the cryptos
# worksheet: it refreshes the Price column (D) for most rows, and for
# rows 42/43 it swaps the BKEXToken/CEJI entries (Coin, Link, Price,
# Volume) because their ranking order changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    # Force the cell to stay a text value (avoid Excel re-interpreting
    # numeric-looking strings as numbers and losing formatting such as
    # trailing zeros).
    $rng.NumberFormat = "@"
    $rng.Value = $Text
}

# --- Price (column D) updates -------------------------------------------
Set-TextValue "D2"  "276.27"
Set-TextValue "D3"  "20.94"
Set-TextValue "D4"  "6.218"
Set-TextValue "D5"  "0.06200"
Set-TextValue "D6"  "3.580"
Set-TextValue "D7"  "6.550"
Set-TextValue "D8"  "1.483"
Set-TextValue "D10" "0.01385"
Set-TextValue "D12" "0.08239"
Set-TextValue "D13" "0.03506"
Set-TextValue "D14" "0.03109"
Set-TextValue "D16" "3.770"
Set-TextValue "D17" "0.001620"
Set-TextValue "D18" "0.04688"
Set-TextValue "D19" "0.006437"
Set-TextValue "D20" "0.006152"
Set-TextValue "D23" "3.823"
Set-TextValue "D24" "2.358"
Set-TextValue "D26" "0.1232"
Set-TextValue "D28" "0.0002739"
Set-TextValue "D40" "0.04674"
Set-TextValue "D41" "0.007026"
Set-TextValue "D44" "0.01089"
Set-TextValue "D45" "0.00006184"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "D47" "0.8459"
Set-TextValue "D48" "0.001969"

# --- Rows 42/43: BKEXToken and CEJI swap places ---------------------------
# Row 42 becomes CEJI, row 43 becomes BKEXToken (A column / rank stays put).
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.004603"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1101"
$ws.Range("E43").Value = "42BKEXTokenBKK"

Write-Host "Applied symbol list update"
